# clear transfer1_ & transfer.py
# Rewrite column K ("transfer") so that it holds the same plain,
# comma-separated text already present in column J ("ref") instead of
# a Python tuple-repr string like ('A', 'B').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $kCell = $ws.Cells.Item($r, 11)   # column K = transfer
    $kVal = $kCell.Value()

    if ($kVal -ne $null -and $kVal -ne "") {
        $jCell = $ws.Cells.Item($r, 10)   # column J = ref
        $jVal = $jCell.Value()
        $kCell.Value = $jVal
    }
}
